# Adds a second commodity ("TABLES") to the DEAS Equipment workbook:
#  - one new row in "Room Inventories" (S1 / TABLES / 10)
#  - a new "TABLES" row in "Commodities", and CHAIRS quantities rearranged
#  - mirrored rows (Equipment Type = TABLES) appended to the four arc sheets
#    (Movement Arcs, Storage Room Arcs, Event Room Arcs, Utility Arcs)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Room Inventories: add row 3 (S1, TABLES, 10)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Room Inventories")
$ws.Cells.Item(3,1).Value = "S1"
$ws.Cells.Item(3,2).Value = "TABLES"
$ws.Cells.Item(3,3).Value = 10
$ws.Range("F11").Select()

# ---------------------------------------------------------------------
# Commodities: CHAIRS becomes 2 units/parcel, 1 sqft volume/parcel;
# add TABLES as 1 unit/parcel, 4 sqft volume/parcel
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Commodities")
$ws.Cells.Item(2,2).Value = 2
$ws.Cells.Item(2,3).Value = 1
$ws.Cells.Item(3,1).Value = "TABLES"
$ws.Cells.Item(3,2).Value = 1
$ws.Cells.Item(3,3).Value = 4
$ws.Range("C2").Select()

# ---------------------------------------------------------------------
# Movement Arcs: append rows 12-21 mirroring rows 2-11 but for TABLES
# (Equipment Type -> TABLES, Uij capacity 72 -> 10)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Movement Arcs")
$rows = @(
    @("E1", 0, "b", "E1", 1, "a", "TABLES", 0, 10, 0),
    @("E1", 0, "b", "S1", 1, "a", "TABLES", 0, 10, 7),
    @("S1", 0, "b", "E1", 1, "a", "TABLES", 0, 10, 7),
    @("S1", 0, "b", "S1", 1, "a", "TABLES", 0, 10, 0),
    @("E1", 1, "b", "E1", 2, "a", "TABLES", 0, 10, 0),
    @("E1", 1, "b", "S1", 2, "a", "TABLES", 0, 10, 7),
    @("S1", 1, "b", "E1", 2, "a", "TABLES", 0, 10, 7),
    @("S1", 1, "b", "S1", 2, "a", "TABLES", 0, 10, 0),
    @("S1", 2, "b", "t",  3, "a", "TABLES", 0, 10, 0),
    @("E1", 2, "b", "t",  3, "a", "TABLES", 0, 10, 0)
)
$r = 12
foreach ($row in $rows) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
$ws.Range("M14").Select()

# ---------------------------------------------------------------------
# Storage Room Arcs: append rows 4-5 for TABLES
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Storage Room Arcs")
$rows = @(
    @("S1", 1, "a", "S1", 1, "b", "TABLES", 0, 10, 0),
    @("S1", 2, "a", "S1", 2, "b", "TABLES", 0, 10, 0)
)
$r = 4
foreach ($row in $rows) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
$ws.Range("I5").Select()

# ---------------------------------------------------------------------
# Utility Arcs: append rows 5-7 for TABLES
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Utility Arcs")
$rows = @(
    @("s", 0, "a", "E1", 0, "b", "TABLES", 6, 6, 0),
    @("s", 0, "a", "S1", 0, "b", "TABLES", 4, 4, 0),
    @("t", 3, "a", "t",  3, "b", "TABLES", 10, 10, 0)
)
$r = 5
foreach ($row in $rows) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
$ws.Range("D7").Select()

# ---------------------------------------------------------------------
# Event Room Arcs: append rows 4-5 for TABLES
# (touched last so it ends up the workbook's active sheet/tab, matching
#  the saved view: firstSheet 3, activeTab 7 -> "Event Room Arcs")
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Event Room Arcs")
$rows = @(
    @("E1", 1, "a", "E1", 1, "b", "TABLES", 3, 3, 0),
    @("E1", 2, "a", "E1", 2, "b", "TABLES", 9, 9, 0)
)
$r = 4
foreach ($row in $rows) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
$ws.Activate()
$ws.Range("G5").Select()

Write-Output "done"
